# Adds a new sheet "Hoja1 (4)" - a copy of "Hoja1 (3)" - summarizing
# problematics ("N") vs knowledge areas, per commit message:
#   "Se agrega tabla con resumen de problematicas vs las areas de conocimiento"

$wb = $excel.ActiveWorkbook

# 1) Duplicate "Hoja1 (3)" and place the copy right after it (becomes the
#    new last sheet / active sheet, matching Excel's "Move or Copy" result).
$wsSource = $wb.Worksheets.Item("Hoja1 (3)")
$wsSource.Copy($null, $wsSource) | Out-Null

$wsNew = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsNew.Name = "Hoja1 (4)"

# 2) On the new sheet, replace every category marker ("Inicialización",
#    "Dirección", "Planificación", "Calidad", "Personas", "Riesgos") with a
#    single generic "N" mark - the new table only flags which
#    problem/process-group intersections apply.
$markedCells = @("C4", "E4", "F4", "C5", "D5", "D6", "E8", "F8", "D9", "E9", "D11", "E11")
foreach ($addr in $markedCells) {
    $wsNew.Range($addr).Value = "N"
}

# 3) Drop the leftover helper column (K4:K9) that isn't part of the new table.
$wsNew.Range("K4:K9").Clear()

# 4) Restore sensible selections: the old sheet keeps a normal selection,
#    while the newly added (now active) sheet gets its own.
$wsNew.Activate()
$wsNew.Range("I6").Select() | Out-Null

$wsSource.Activate()
$wsSource.Range("C4").Select() | Out-Null

$wsNew.Activate()
